$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Monitorista"
$ws.Range("B2").Value = "SISTEMAS DE SEGURIDAD PRIVADA ARGOS S.A. DE C.V."
$ws.Range("C2").Value = "Torreón, Coah."
$ws.Range("D2").Value = "`$10,000 por mes"
$ws.Range("E2").Value = "Postulación vía Indeed"
$ws.Range("F2").Value = "Indeed-Programador"

$ws.Range("A3").Value = "Desarrollador Web"
$ws.Range("B3").Value = "SISTEMAS DE SEGURIDAD PRIVADA ARGOS S.A. DE C.V."
$ws.Range("C3").Value = "Torreón, Coah."
$ws.Range("D3").Value = "`$13,000 a `$15,000 por mes"
$ws.Range("E3").Value = "Postulación vía Indeed"
$ws.Range("F3").Value = "Indeed-Programador"

$ws.Range("A4").Value = "ESPECILISTA DENTAL para trabajar en Torreón ( ingreso mínimo asegurado)"
$ws.Range("B4").Value = "LOS DENTISTAS"
$ws.Range("C4").Value = "Torreón, Coah."
$ws.Range("D4").Value = "`$25,000 a `$50,000 por mes"
$ws.Range("E4").Value = "Postulación vía Indeed"
$ws.Range("F4").Value = "Indeed-Programador"

$ws.Range("A5").Value = "Supervisor de Atención a Clientes"
$ws.Range("B5").Value = "Paquetexpress"
$ws.Range("C5").Value = "Torreón, Coah."
$ws.Range("D5").Value = "`$14,900 por mes"
$ws.Range("E5").Value = "Postulación vía Indeed"
$ws.Range("F5").Value = "Indeed-Programador"

$ws.Range("A6").Value = "ANALISTA GASTOS DE VIAJE"
$ws.Range("B6").Value = "Grupo LALA"
$ws.Range("C6").Value = "Torreón, Coah."
$ws.Range("D6").Value = "N/A"
$ws.Range("E6").Value = "N/A"
$ws.Range("F6").Value = "Indeed-Programador"

$ws.Range("A7").Value = "ASESOR DE SERVICIO AUTOMOTRIZ"
$ws.Range("B7").Value = "GRUPO ALAMEDA"
$ws.Range("C7").Value = "Torreón, Coah."
$ws.Range("D7").Value = "`$2,000 a `$3,500 por semana"
$ws.Range("E7").Value = "Postulación vía Indeed"
$ws.Range("F7").Value = "Indeed-Programador"

$ws.Range("A8").Value = "MESERO/VENDEDOR"
$ws.Range("B8").Value = "PUCCINO´S - Restaurante-Bar Italiano"
$ws.Range("C8").Value = "Torreón, Coah."
$ws.Range("D8").Value = "Tiempo completo"
$ws.Range("E8").Value = "Postulación vía Indeed"
$ws.Range("F8").Value = "Indeed-Programador"

$ws.Range("A9").Value = "Auxiliar Administrativo"
$ws.Range("B9").Value = "Constructora Cayarga, S.A de C.V."
$ws.Range("C9").Value = "27110, Residencial los Llanos, Coah."
$ws.Range("D9").Value = "N/A"
$ws.Range("E9").Value = "N/A"
$ws.Range("F9").Value = "Indeed-Programador"

$ws.Range("A10").Value = "Recepcionista Español - Dirección Comercial"
$ws.Range("B10").Value = "Peñoles"
$ws.Range("C10").Value = "Torreón, Coah."
$ws.Range("D10").Value = "`$10,000 a `$12,000 por mes"
$ws.Range("E10").Value = "Postulación vía Indeed"
$ws.Range("F10").Value = "Indeed-Programador"

$ws.Range("A11").Value = "CAJERO (A)"
$ws.Range("B11").Value = "SEPSA"
$ws.Range("C11").Value = "27000, Torreón Centro, Coah."
$ws.Range("D11").Value = "N/A"
$ws.Range("E11").Value = "N/A"
$ws.Range("F11").Value = "Indeed-Programador"

$ws.Range("A12").Value = "SUPERVISOR ALMACEN"
$ws.Range("B12").Value = "Grupo Modelo"
$ws.Range("C12").Value = "Torreón, Coah."
$ws.Range("D12").Value = "N/A"
$ws.Range("E12").Value = "N/A"
$ws.Range("F12").Value = "Indeed-Programador"

$ws.Range("A13").Value = "Asistente Administrativo"
$ws.Range("B13").Value = "Easy Way Products"
$ws.Range("C13").Value = "27280, Eduardo Guerra, Coah."
$ws.Range("D13").Value = "Tiempo completo"
$ws.Range("E13").Value = "Postulación vía Indeed"
$ws.Range("F13").Value = "Indeed-Programador"

$ws.Range("A14").Value = "RECEPCIONISTA"
$ws.Range("B14").Value = "GAFI SERVICIOS SA DE CV"
$ws.Range("C14").Value = "Torreón, Coah."
$ws.Range("D14").Value = "N/A"
$ws.Range("E14").Value = "N/A"
$ws.Range("F14").Value = "Indeed-Programador"

$ws.Range("A15").Value = "Manufacturing Engr Manager"
$ws.Range("B15").Value = "Caterpillar"
$ws.Range("C15").Value = "Torreón, Coah."
$ws.Range("D15").Value = "Tiempo completo"
$ws.Range("E15").Value = "N/A"
$ws.Range("F15").Value = "Indeed-Programador"

$ws.Range("A16").Value = "Supervisor, Production"
$ws.Range("B16").Value = "YAZAKI Corporation"
$ws.Range("C16").Value = "Torreón, Coah."
$ws.Range("D16").Value = "N/A"
$ws.Range("E16").Value = "N/A"
$ws.Range("F16").Value = "Indeed-Programador"

